$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows above row 49 (old rows 49-53 were blank spacer rows;
# the new layout needs 8 "cycle_threading" variable rows (49-56) plus
# 3 blank spacer rows (57-59) before the 5001.. block, i.e. 6 more rows
# than before). Inserting here shifts everything below (5001 block,
# 5380/5398 block, turning-speed calc block) down by 6 rows, matching
# the target layout (old row 54 -> new row 60, old row 81 -> new row 87).
$ws.Rows("49:54").Insert()

# New shared string used by the inserted rows.
$cycleName = "cycle_threading"

$firstVal = 1550
for ($i = 0; $i -lt 8; $i++) {
    $r = 49 + $i
    $val = $firstVal + $i

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $val
    $aCell.Interior.Color = 65535     # yellow, matches the other 14xx/15xx rows
    $aCell.Borders.LineStyle = -4142  # xlLineStyleNone

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value = $cycleName
    $cCell.Interior.Pattern = -4142   # xlNone
    $cCell.Borders.LineStyle = -4142  # xlLineStyleNone

    # Rows 49-54 also carry an (empty) formatted B cell, like the C
    # column, mirroring the 1500-1512 block above. Rows 55-56 don't get
    # a B cell at all.
    if ($i -lt 6) {
        $bCell = $ws.Cells.Item($r, 2)
        $bCell.Interior.Pattern = -4142   # xlNone
        $bCell.Borders.LineStyle = -4142  # xlLineStyleNone
    }
}

# Reposition the view roughly where the author left it (best effort —
# the window/scroll chrome itself is session state, but the selected
# cell is part of the saved sheet view).
$ws.Range("D53").Select()
